# Apply edits to "Final Checklist" sheet (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in new "Your Points" values (column F) for several rows
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("F9").Value = 0.5
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 0.5

# Add the SUM formula to the total row for column F
$ws.Range("F14").Formula = "=SUM(F3:F13)"

# Update the active selection to match the final state of the file
$ws.Range("F15").Select()
